# Weekly update: insert a new Mandarina price record at row 102.
# All existing records from row 102 down to row 114 shift down by one
# row (to 103-115), and the new record's data is written into row 102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 102, pushing existing rows 102:114
# down to 103:115 (and extending the used range to row 115).
$ws.Rows("102:102").Insert()

# Populate the newly inserted row 102 with this week's data.
$ws.Cells.Item(102, 1).Value = 1
$ws.Cells.Item(102, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(102, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(102, 4).Value = 44783
$ws.Cells.Item(102, 5).Value = 15
$ws.Cells.Item(102, 6).Value = "Fruta"
$ws.Cells.Item(102, 7).Value = 100102
$ws.Cells.Item(102, 8).Value = "Cítricos"
$ws.Cells.Item(102, 9).Value = 100102004
$ws.Cells.Item(102, 10).Value = "Mandarina"
$ws.Cells.Item(102, 11).Value = "Clemenuless"
$ws.Cells.Item(102, 12).Value = "Segunda"
$ws.Cells.Item(102, 13).Value = 300
$ws.Cells.Item(102, 14).Value = 15000
$ws.Cells.Item(102, 15).Value = 16000
$ws.Cells.Item(102, 16).Value = 15500
$ws.Cells.Item(102, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(102, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(102, 19).Value = 775
$ws.Cells.Item(102, 20).Value = 20
